{"js": "// Change the time stamp at the end of the document from \"15:56\" to\n// \"16:05\" (commit message: \"Modifica delle ore 16:05\").\n//\n// The last paragraph in the body is: <bookmarkStart \"_GoBack\"/><bookmarkEnd/>\n// followed by a run containing \"15:56\". In the target revision the run's\n// text becomes \"16:05\" and the (still-empty) \"_GoBack\" bookmark is moved so\n// it again sits right after the edited text - i.e. exactly what Word does\n// when you select the old text and type the replacement: the \"_GoBack\"\n// bookmark is re-anchored to the point of the most recent edit.\n\nconst body = context.document.body;\n\n// 1) Find the run that holds the old time and swap its text in place so\n//    every other run/paragraph property is left untouched.\nconst oldTime = body.search(\"15:56\", { matchCase: true, matchWholeWord: false });\noldTime.load(\"text\");\nawait context.sync();\n\nif (oldTime.items.length === 0) {\n  throw new Error('Could not find \"15:56\" in the document body.');\n}\n\nconst oldRange = oldTime.items[0];\noldRange.insertText(\"16:05\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Re-anchor the \"_GoBack\" bookmark so it sits after the new text again\n//    (matching Word's own behaviour of tracking the last edit point).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst newTime = body.search(\"16:05\", { matchCase: true, matchWholeWord: false });\nnewTime.load(\"text\");\nawait context.sync();\n\nconst newRange = newTime.items[0];\nconst afterNewText = newRange.getRange(Word.RangeLocation.end);\nafterNewText.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Change the time stamp at the end of the document from \"15:56\" to\n# \"16:05\" (commit message: \"Modifica delle ore 16:05\").\n#\n# The last paragraph in the body is:\n#   <w:bookmarkStart w:name=\"_GoBack\"/><w:bookmarkEnd/><w:r>...15:56...</w:r>\n# In the target revision the run's text becomes \"16:05\" and the (still\n# empty) \"_GoBack\" bookmark is moved so it sits right after the edited\n# text instead of before it - exactly what Word does when you select the\n# old text and type a replacement: \"_GoBack\" is re-anchored to the point\n# of the most recent edit.\n\n$d = $word.ActiveDocument\n\n# 1) Replace the time text in place so every other run/paragraph property\n#    is left untouched.\n$findRange = $d.Content\n$find = $findRange.Find\n$find.Text = \"15:56\"\n$find.Replacement.Text = \"16:05\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n# 2) Drop the old \"_GoBack\" bookmark - it will be re-added after the new text.\n$d.Bookmarks(\"_GoBack\").Delete()\n\n# 3) Find the freshly inserted text and collapse a range to right after it.\n$locateRange = $d.Content\n$find2 = $locateRange.Find\n$find2.Text = \"16:05\"\n$find2.Execute() | Out-Null\n$locateRange.Collapse(0)   # wdCollapseEnd\n\n# A collapsed range sitting immediately before a paragraph mark confuses\n# Bookmarks.Add, so temporarily insert a throwaway character after the\n# collapse point, anchor the bookmark next to ordinary text, then remove\n# the throwaway character again.\n$locateRange.InsertAfter(\"Z\")\n\n$bmRange = $d.Range($locateRange.Start, $locateRange.Start)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange) | Out-Null\n\n$placeholder = $d.Range($bmRange.Start, $bmRange.Start + 1)\n$placeholder.Delete() | Out-Null\n"}
